$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 2208
$ws.Range("I8").Value = 2208
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 6624
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -6485
$ws.Range("N8").ClearContents()

$ws.Range("H15").Value = 271.14
$ws.Range("I15").Value = 271.14
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 813.42
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -644.42

$ws.Range("H33").Value = 479.35
$ws.Range("I33").Value = 508.16666
$ws.Range("J33").Value = 220
$ws.Range("K33").Value = 508.16666
$ws.Range("L33").Value = 220
$ws.Range("M33").Value = -279.16666
$ws.Range("N33").Value = -678

$ws.Range("H64").Value = 2973.913
$ws.Range("I64").Value = 2945.4546
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 2945.4546
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2697.4546
$ws.Range("N64").Value = -3496

$ws.Range("H67").Value = 2973.913
$ws.Range("I67").Value = 2945.4546
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 2945.4546
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2087.4546
$ws.Range("N67").Value = -4716

$ws.Range("H100").Value = 4235.6665
$ws.Range("I100").Value = 2135
$ws.Range("J100").Value = 6336.3335
$ws.Range("K100").Value = 2135
$ws.Range("L100").Value = 6336.3335
$ws.Range("M100").Value = -1594
$ws.Range("N100").Value = -7418.3335

$ws.Range("H107").Value = 2387
$ws.Range("I107").Value = 2560.375
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 2560.375
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = -640.375
$ws.Range("N107").Value = -4840

$ws.Range("H113").Value = 2382.85
$ws.Range("I113").Value = 2017
$ws.Range("J113").Value = 2830
$ws.Range("K113").Value = 2017
$ws.Range("L113").Value = 2830
$ws.Range("M113").Value = 1237
$ws.Range("N113").Value = -9338

$ws.Range("H115").Value = 2128.5
$ws.Range("I115").Value = 1547.5
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 4642.5
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -3075.5
$ws.Range("N115").Value = -12134

$ws.Range("H127").Value = 1230.3077
$ws.Range("I127").Value = 500
$ws.Range("J127").Value = 1686.75
$ws.Range("K127").Value = 1500
$ws.Range("L127").Value = 5060.25
$ws.Range("M127").Value = 3460
$ws.Range("N127").Value = -14980.25

$ws.Range("H137").Value = 1667.409
$ws.Range("I137").Value = 1327.3636
$ws.Range("J137").Value = 2007.4546
$ws.Range("K137").Value = 3982.0908
$ws.Range("L137").Value = 6022.3638
$ws.Range("M137").Value = -1432.0908
$ws.Range("N137").Value = -11122.3638

$ws.Range("H138").Value = 3506.195
$ws.Range("I138").Value = 3244.0667
$ws.Range("J138").Value = 3564.8806
$ws.Range("K138").Value = 9732.2001
$ws.Range("L138").Value = 10694.6418
$ws.Range("M138").Value = -4592.2001
$ws.Range("N138").Value = -20974.6418

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1904370
$ws.Range("I32").Value = 2355260.2
$ws.Range("J32").Value = 21240.059
$ws.Range("K32").Value = 2355260.2
$ws.Range("L32").Value = 21240.059
$ws.Range("M32").Value = -2354973.2
$ws.Range("N32").Value = -21814.059

$ws.Range("H74").Value = 2710.7097
$ws.Range("I74").Value = 1908.6364
$ws.Range("J74").Value = 4671.3335
$ws.Range("K74").Value = 1908.6364
$ws.Range("L74").Value = 4671.3335
$ws.Range("M74").Value = -1034.6364
$ws.Range("N74").Value = -6419.3335

$ws.Range("H77").Value = 2710.7097
$ws.Range("I77").Value = 1908.6364
$ws.Range("J77").Value = 4671.3335
$ws.Range("K77").Value = 9543.182000000001
$ws.Range("L77").Value = 23356.6675
$ws.Range("M77").Value = -5175.182000000001
$ws.Range("N77").Value = -32092.6675

$ws.Range("H102").Value = 3222.7144
$ws.Range("I102").Value = 3280
$ws.Range("J102").Value = 3199.8
$ws.Range("K102").Value = 3280
$ws.Range("L102").Value = 3199.8
$ws.Range("M102").Value = -1658
$ws.Range("N102").Value = -6443.8

$ws.Range("H122").Value = 57133.332
$ws.Range("I122").Value = 68173.336
$ws.Range("J122").Value = 1933.3334
$ws.Range("K122").Value = 204520.008
$ws.Range("L122").Value = 5800.0002
$ws.Range("M122").Value = -202070.008
$ws.Range("N122").Value = -10700.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 864
$ws.Range("I94").Value = 790
$ws.Range("J94").Value = 975
$ws.Range("K94").Value = 790
$ws.Range("L94").Value = 975
$ws.Range("M94").Value = -339
$ws.Range("N94").Value = -1877

$ws.Range("H134").Value = 2767.4866
$ws.Range("I134").Value = 2718.0908
$ws.Range("J134").Value = 3175
$ws.Range("K134").Value = 8154.2724
$ws.Range("L134").Value = 9525
$ws.Range("M134").Value = -5619.2724
$ws.Range("N134").Value = -14595

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 805.0526
$ws.Range("I105").Value = 765
$ws.Range("J105").Value = 891.8333
$ws.Range("K105").Value = 765
$ws.Range("L105").Value = 891.8333
$ws.Range("M105").Value = 982
$ws.Range("N105").Value = -4385.8333

$ws.Range("H132").Value = 5378694
$ws.Range("I132").Value = 2166.647
$ws.Range("J132").Value = 11907334
$ws.Range("K132").Value = 6499.941
$ws.Range("L132").Value = 35722002
$ws.Range("M132").Value = -3969.941
$ws.Range("N132").Value = -35727062

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 648.4
$ws.Range("I7").Value = 374.44446
$ws.Range("J7").Value = 765.8095
$ws.Range("K7").Value = 1123.33338
$ws.Range("L7").Value = 2297.4285
$ws.Range("M7").Value = -1011.33338
$ws.Range("N7").Value = -2521.4285

$ws.Range("H86").Value = 849.5
$ws.Range("I86").Value = 849.5
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2548.5
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -1362.5
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 849.5
$ws.Range("I89").Value = 849.5
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 7645.5
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -1717.5
$ws.Range("N89").ClearContents()

$ws.Range("H134").Value = 6563.641
$ws.Range("I134").Value = 3349.9
$ws.Range("J134").Value = 7671.8276
$ws.Range("K134").Value = 10049.7
$ws.Range("L134").Value = 23015.4828
$ws.Range("M134").Value = -4979.700000000001
$ws.Range("N134").Value = -33155.4828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H118").Value = 43935
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 43935
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 43935
$ws.Range("N118").Value = -47249

$ws.Range("H121").Value = 36987.25
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 36987.25
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 36987.25
$ws.Range("N121").Value = -40481.25

$ws.Range("H122").Value = 7617.8237
$ws.Range("I122").Value = 10672.363
$ws.Range("J122").Value = 2017.8334
$ws.Range("K122").Value = 32017.089
$ws.Range("L122").Value = 6053.5002
$ws.Range("M122").Value = -29567.089
$ws.Range("N122").Value = -10953.5002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12550.467
$ws.Range("I22").Value = 612.375
$ws.Range("J22").Value = 26194
$ws.Range("K22").Value = 612.375
$ws.Range("L22").Value = 26194
$ws.Range("M22").Value = -317.375
$ws.Range("N22").Value = -26784

$ws.Range("H27").Value = 12550.467
$ws.Range("I27").Value = 612.375
$ws.Range("J27").Value = 26194
$ws.Range("K27").Value = 612.375
$ws.Range("L27").Value = 26194
$ws.Range("M27").Value = -505.375
$ws.Range("N27").Value = -26408

$ws.Range("H61").Value = 3312.875
$ws.Range("I61").Value = 2339.2144
$ws.Range("J61").Value = 4676
$ws.Range("K61").Value = 2339.2144
$ws.Range("L61").Value = 4676
$ws.Range("M61").Value = -2137.2144
$ws.Range("N61").Value = -5080

$ws.Range("H100").Value = 3698
$ws.Range("I100").Value = 3240
$ws.Range("J100").Value = 4614
$ws.Range("K100").Value = 3240
$ws.Range("L100").Value = 4614
$ws.Range("M100").Value = -2699
$ws.Range("N100").Value = -5696

$ws.Range("H113").Value = 3312.875
$ws.Range("I113").Value = 2339.2144
$ws.Range("J113").Value = 4676
$ws.Range("K113").Value = 2339.2144
$ws.Range("L113").Value = 4676
$ws.Range("M113").Value = -169.2143999999998
$ws.Range("N113").Value = -9016

$ws.Range("H122").Value = 4523.8096
$ws.Range("I122").Value = 4060
$ws.Range("J122").Value = 4668.75
$ws.Range("K122").Value = 12180
$ws.Range("L122").Value = 14006.25
$ws.Range("M122").Value = -9730
$ws.Range("N122").Value = -18906.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 965.7727
$ws.Range("I113").Value = 1122.1765
$ws.Range("J113").Value = 434
$ws.Range("K113").Value = 3366.5295
$ws.Range("L113").Value = 1302
$ws.Range("M113").Value = -1196.5295
$ws.Range("N113").Value = -5642
